$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line to add "EA (8), " before "EP (8)"
$d.Content.Find.Execute("Curso (semestre ideal): EP (8)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EA (8), EP (8)", 2)

# 2. Remove the trailing "Requisitos" heading paragraph and the
#    "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)" bullet paragraph
#    that follows it (the last two paragraphs of the document).
$count = $d.Paragraphs.Count
$reqHeading = $d.Paragraphs.Item($count - 1)
$reqBullet = $d.Paragraphs.Item($count)
$r = $d.Range($reqHeading.Range.Start, $reqBullet.Range.End)
$r.Delete()
